$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 48,11
$data[0,0] = ""
$data[0,1] = "Q0"
$data[0,2] = "Q1"
$data[0,3] = "Q2"
$data[0,4] = "Q3"
$data[0,5] = "Q4"
$data[0,6] = "Q5"
$data[0,7] = "Q6"
$data[0,8] = "Q7"
$data[0,9] = "Q8"
$data[0,10] = "Q9"
$data[1,0] = "2010-04-01 00:00:00_diff"
$data[1,1] = -2.599775093864654
$data[1,2] = 0.4213078533922215
$data[1,3] = -0.0851732134718878
$data[1,4] = -0.5545455541496284
$data[1,5] = -0.2341436367370153
$data[1,6] = -0.02969014145458965
$data[1,7] = -0.1524194508597342
$data[1,8] = 0.5892983349069034
$data[1,9] = 0.5410791075487134
$data[1,10] = 0.2469615835596112
$data[2,0] = "2010-10-01 00:00:00_diff"
$data[2,1] = 0.7599138577479545
$data[2,2] = -0.8296272591627797
$data[2,3] = 0.183025357444689
$data[2,4] = -0.1193166304426081
$data[2,5] = 0.1037808198054395
$data[2,6] = 0.5999655385994425
$data[2,7] = 0.724112664506922
$data[2,8] = 0.3076557451057896
$data[2,9] = 0.3331284835103464
$data[2,10] = 0.3012303191900245
$data[3,0] = "2011-04-01 00:00:00_diff"
$data[3,1] = -0.1932672371045038
$data[3,2] = 0.04324661062325436
$data[3,3] = -0.01161065538656181
$data[3,4] = 0.6246649691719768
$data[3,5] = 0.6777283932422284
$data[3,6] = 0.297938859626624
$data[3,7] = 0.3044483231023983
$data[3,8] = 0.2825089781401424
$data[3,9] = -0.0973857149104142
$data[3,10] = -0.4061154638499449
$data[4,0] = "2011-10-01 00:00:00_diff"
$data[4,1] = 0.01885790495154901
$data[4,2] = 0.6365171369896647
$data[4,3] = 0.68016851135583
$data[4,4] = 0.31279625129615
$data[4,5] = 0.3063033874885817
$data[4,6] = 0.2942571669794643
$data[4,7] = -0.0931641912349494
$data[4,8] = -0.3966786672555065
$data[4,9] = 0.6137917488812377
$data[4,10] = -1.55398144811108
$data[5,0] = "2012-04-01 00:00:00_diff"
$data[5,1] = 0.7680233526332453
$data[5,2] = 0.086658463984111
$data[5,3] = 0.3550562147590637
$data[5,4] = 0.1996620052474151
$data[5,5] = -0.0792244915534977
$data[5,6] = -0.4452253649762511
$data[5,7] = 0.6091242210487386
$data[5,8] = -1.585289023135289
$data[5,9] = 0.4557940629208951
$data[5,10] = 0.4460097931269203
$data[6,0] = "2012-10-01 00:00:00_diff"
$data[6,1] = 0.1149817136733693
$data[6,2] = 0.1886961826570769
$data[6,3] = -0.1743510322417958
$data[6,4] = -0.4576880558317372
$data[6,5] = 0.5575615649551221
$data[6,6] = -1.605488280710164
$data[6,7] = 0.4185686116426709
$data[6,8] = 0.4210803793827069
$data[6,9] = -0.8032415029445406
$data[6,10] = -0.2573748024414758
$data[7,0] = "2013-04-01 00:00:00_diff"
$data[7,1] = -0.1984711623028061
$data[7,2] = -0.522343575111026
$data[7,3] = 0.5457778617256782
$data[7,4] = -1.635245477730326
$data[7,5] = 0.4063622901268292
$data[7,6] = 0.4015824023328216
$data[7,7] = -0.8167612454749933
$data[7,8] = -0.2737209468961671
$data[7,9] = 0.4522297872112956
$data[7,10] = -0.49062832164035
$data[8,0] = "2013-10-01 00:00:00_diff"
$data[8,1] = 0.578229357260996
$data[8,2] = -1.451239549016748
$data[8,3] = 0.4322360342429518
$data[8,4] = 0.4897976924316668
$data[8,5] = -0.7812097065401281
$data[8,6] = -0.2140188449768187
$data[8,7] = 0.4939837660135075
$data[8,8] = -0.4398144183660089
$data[8,9] = -0.5732498698339226
$data[8,10] = -0.5255200426117579
$data[9,0] = "2014-04-01 00:00:00_diff"
$data[9,1] = 0.03913348239479289
$data[9,2] = 1.027270011939213
$data[9,3] = -0.9662182088689218
$data[9,4] = 0.006854978822963043
$data[9,5] = 0.4253527326840144
$data[9,6] = -0.3358070042503943
$data[9,7] = -0.5867795836918103
$data[9,8] = -0.4666029924297193
$data[9,9] = 0.4288730018125405
$data[9,10] = -0.1606097620349077
$data[10,0] = "2014-10-01 00:00:00_diff"
$data[10,1] = -0.8655550547631281
$data[10,2] = -0.3643995038402464
$data[10,3] = 0.460959376748889
$data[10,4] = -0.5052523480449113
$data[10,5] = -0.6030493650709703
$data[10,6] = -0.5684950311133089
$data[10,7] = 0.3862892789456705
$data[10,8] = -0.238319052841212
$data[10,9] = -0.4163595244981116
$data[10,10] = 0.0805329408679176
$data[11,0] = "2015-04-01 00:00:00_diff"
$data[11,1] = 0.6481750507136302
$data[11,2] = -0.4116222461567253
$data[11,3] = -0.5414319719249838
$data[11,4] = -0.5241108785450015
$data[11,5] = 0.4246184323196733
$data[11,6] = -0.2030005133752568
$data[11,7] = -0.3823142582379787
$data[11,8] = 0.1141332747280269
$data[11,9] = -0.8852847537018906
$data[11,10] = -0.152826881537631
$data[12,0] = "2015-10-01 00:00:00_diff"
$data[12,1] = -0.7565604878827148
$data[12,2] = -0.4057663831855379
$data[12,3] = 0.3522497996303807
$data[12,4] = -0.1715355603405611
$data[12,5] = -0.4092831094845341
$data[12,6] = 0.1193314088729555
$data[12,7] = -0.8980522263979865
$data[12,8] = -0.1556535701336084
$data[12,9] = 0.1348417040676331
$data[12,10] = -0.1120512842982363
$data[13,0] = "2016-04-01 00:00:00_diff"
$data[13,1] = 0.4998011962756446
$data[13,2] = -0.0734274294216598
$data[13,3] = -0.3598965770760573
$data[13,4] = 0.1625872521950072
$data[13,5] = -0.864103690938309
$data[13,6] = -0.1222990917425359
$data[13,7] = 0.1664338649438454
$data[13,8] = -0.08049589489028602
$data[13,9] = -0.5730541687050947
$data[13,10] = -0.4958239707762366
$data[14,0] = "2016-10-01 00:00:00_diff"
$data[14,1] = -0.4908112086026729
$data[14,2] = 0.1874428494022752
$data[14,3] = -0.9058127297462113
$data[14,4] = -0.1202113871632399
$data[14,5] = 0.1477633332452291
$data[14,6] = -0.0865810009705153
$data[14,7] = -0.5854697566684117
$data[14,8] = -0.5045745169269311
$data[14,9] = 0.2842444055692547
$data[14,10] = 0.094314310433458
$data[15,0] = "2017-04-01 00:00:00_diff"
$data[15,1] = -0.7610382668704541
$data[15,2] = -0.1697185737259005
$data[15,3] = 0.1975744017561182
$data[15,4] = -0.09296950593211634
$data[15,5] = -0.5624135302179352
$data[15,6] = -0.4978696587057602
$data[15,7] = 0.2996381065732735
$data[15,8] = 0.1049343773025755
$data[15,9] = -0.1810904405051886
$data[15,10] = 0.7100950095189205
$data[16,0] = "2017-10-01 00:00:00_diff"
$data[16,1] = 0.3853124329624699
$data[16,2] = -0.0636353541259792
$data[16,3] = -0.5057317212606856
$data[16,4] = -0.4866080797737675
$data[16,5] = 0.3244420063798164
$data[16,6] = 0.1159782022070233
$data[16,7] = -0.1646040460938459
$data[16,8] = 0.7222453785596309
$data[16,9] = -0.3978929871876867
$data[16,10] = 1.818769300392457
$data[17,0] = "2018-01-01 00:00:00_diff"
$data[17,1] = -0.0411574570261935
$data[17,2] = -0.6196320408848983
$data[17,3] = -0.4817670185010929
$data[17,4] = 0.2802806761101042
$data[17,5] = 0.1075973650542067
$data[17,6] = -0.1897256205557853
$data[17,7] = 0.7081629777537968
$data[17,8] = -0.4175472274164846
$data[17,9] = 1.802574501419392
$data[17,10] = 7.532660084419956
$data[18,0] = "2018-04-01 00:00:00_diff"
$data[18,1] = -0.6322123074925011
$data[18,2] = -0.4742381362309798
$data[18,3] = 0.2721415675329692
$data[18,4] = 0.106451857864849
$data[18,5] = -0.1955958842297873
$data[18,6] = 0.7046241218273783
$data[18,7] = -0.4225415964970449
$data[18,8] = 1.798340658913301
$data[18,9] = 7.527972027331272
$data[18,10] = -4.36761069155567
$data[19,0] = "2018-07-01 00:00:00_diff"
$data[19,1] = -0.5544345664476722
$data[19,2] = 0.4223997105967412
$data[19,3] = 0.06519917072402209
$data[19,4] = -0.1579346443787821
$data[19,5] = 0.6883297515268441
$data[19,6] = -0.4134903260564207
$data[19,7] = 1.791734611061779
$data[19,8] = 7.529266551606639
$data[19,9] = -4.370940361648667
$data[19,10] = 0.5636000270773455
$data[20,0] = "2018-10-01 00:00:00_diff"
$data[20,1] = 0.3026622004506922
$data[20,2] = 0.2048262499793671
$data[20,3] = -0.2046733599345767
$data[20,4] = 0.7259246849823544
$data[20,5] = -0.4292579831143313
$data[20,6] = 1.802332902229133
$data[20,7] = 7.524261510995123
$data[20,8] = -4.367859906123281
$data[20,9] = 0.5620631137850272
$data[20,10] = 1.565731669059103
$data[21,0] = "2019-01-01 00:00:00_diff"
$data[21,1] = 0.2773214709040386
$data[21,2] = -0.2682241111637491
$data[21,3] = 0.7792998334868486
$data[21,4] = -0.4306235619234042
$data[21,5] = 1.839468075212102
$data[21,6] = 7.541092445529947
$data[21,7] = -4.338005966452874
$data[21,8] = 0.5846543025991513
$data[21,9] = 1.592778018715411
$data[21,10] = -2.030234130149802
$data[22,0] = "2019-04-01 00:00:00_diff"
$data[22,1] = -0.2672264875150854
$data[22,2] = 0.6742776619777973
$data[22,3] = -0.4472833318113997
$data[22,4] = 1.783566999271859
$data[22,5] = 7.513849188289668
$data[22,6] = -4.379275238411224
$data[22,7] = 0.5528164771854336
$data[22,8] = 1.556029015863666
$data[22,9] = -2.063849124089159
$data[22,10] = -1.744615694249242
$data[23,0] = "2019-07-01 00:00:00_diff"
$data[23,1] = 0.6086405051256832
$data[23,2] = -0.3644232557353271
$data[23,3] = 1.781188140455046
$data[23,4] = 7.550279024990743
$data[23,5] = -4.364742980045667
$data[23,6] = 0.5774991241904746
$data[23,7] = 1.575067216788767
$data[23,8] = -2.04215268226752
$data[23,9] = -1.724379445304747
$data[23,10] = 0.8571285904013629
$data[24,0] = "2019-10-01 00:00:00_diff"
$data[24,1] = -0.3948898833243785
$data[24,2] = 1.628481131420173
$data[24,3] = 7.529214656055699
$data[24,4] = -4.417393393403793
$data[24,5] = 0.5549235348850442
$data[24,6] = 1.544487103311906
$data[24,7] = -2.065842574713278
$data[24,8] = -1.750067707842675
$data[24,9] = 0.8330224288541155
$data[24,10] = -1.492742449122803
$data[25,0] = "2020-01-01 00:00:00_diff"
$data[25,1] = 1.648811223744012
$data[25,2] = 7.622581430282057
$data[25,3] = -4.39185494626253
$data[25,4] = 0.5915173482346796
$data[25,5] = 1.569371227629475
$data[25,6] = -2.039684916748001
$data[25,7] = -1.725772450408734
$data[25,8] = 0.8573295370653448
$data[25,9] = -1.46868528262261
$data[25,10] = -0.4129619237297894
$data[26,0] = "2020-04-01 00:00:00_diff"
$data[26,1] = 7.430105491556078
$data[26,2] = -4.804399299783788
$data[26,3] = 0.4790604121298505
$data[26,4] = 1.418682353025873
$data[26,5] = -2.128355047819658
$data[26,6] = -1.820804677822804
$data[26,7] = 0.7751543734256706
$data[26,8] = -1.551847621592696
$data[26,9] = -0.4934503117984022
$data[26,10] = -0.275513610932578
$data[27,0] = "2020-07-01 00:00:00_diff"
$data[27,1] = -8.951727421306318
$data[27,2] = -3.240964595099326
$data[27,3] = -1.373514706630456
$data[27,4] = -4.397355256849238
$data[27,5] = -3.619713032289326
$data[27,6] = -0.6960666684218184
$data[27,7] = -2.762904832993337
$data[27,8] = -1.513651065911071
$data[27,9] = -1.149486620479212
$data[27,10] = 0.05092184080034523
$data[28,0] = "2020-10-01 00:00:00_diff"
$data[28,1] = -0.4824776876932789
$data[28,2] = 1.599139081931884
$data[28,3] = -2.176964645512088
$data[28,4] = -1.840977140162832
$data[28,5] = 0.7388118808618853
$data[28,6] = -1.58621144197159
$data[28,7] = -0.5309159314665018
$data[28,8] = -0.3126334443148695
$data[28,9] = 0.7781004828193047
$data[28,10] = 1.194823674284569
$data[29,0] = "2021-01-01 00:00:00_diff"
$data[29,1] = 1.224587888766054
$data[29,2] = -2.056444242200622
$data[29,3] = -1.828588292556533
$data[29,4] = 0.7459962015599261
$data[29,5] = -1.570538489117364
$data[29,6] = -0.516632096391023
$data[29,7] = -0.2986983468163306
$data[29,8] = 0.7922168876341849
$data[29,9] = 1.208927932326396
$data[29,10] = -0.375866296109733
$data[30,0] = "2021-04-01 00:00:00_diff"
$data[30,1] = -1.649879923112767
$data[30,2] = -1.741771509043598
$data[30,3] = 0.6741899500733747
$data[30,4] = -1.543777991623425
$data[30,5] = -0.4957364454691704
$data[30,6] = -0.2907775822781197
$data[30,7] = 0.8048567478956986
$data[30,8] = 1.222081951210451
$data[30,9] = -0.3635527127612387
$data[30,10] = 0.583839642026882
$data[31,0] = "2021-07-01 00:00:00_diff"
$data[31,1] = -2.249670444931811
$data[31,2] = 0.6320633131603571
$data[31,3] = -1.350593266819788
$data[31,4] = -0.4683915811657578
$data[31,5] = -0.248888267974934
$data[31,6] = 0.8716176203667728
$data[31,7] = 1.277643478744205
$data[31,8] = -0.3085675236799043
$data[31,9] = 0.6410202603145777
$data[31,10] = 0.5607407784345422
$data[32,0] = "2021-10-01 00:00:00_diff"
$data[32,1] = -0.1259214769749655
$data[32,2] = -1.475618685980302
$data[32,3] = -0.2170938384811635
$data[32,4] = -0.2416929448297526
$data[32,5] = 0.8856674530228952
$data[32,6] = 1.336194812205091
$data[32,7] = -0.2657147537689935
$data[32,8] = 0.6804779037697479
$data[32,9] = 0.6042710092048686
$data[32,10] = 0.3148496909019965
$data[33,0] = "2022-01-01 00:00:00_diff"
$data[33,1] = -1.456316997758197
$data[33,2] = -0.2221655067966005
$data[33,3] = -0.2298847203035863
$data[33,4] = 0.8954573244745014
$data[33,5] = 1.345928955629549
$data[33,6] = -0.2560777933532004
$data[33,7] = 0.6898230454056837
$data[33,8] = 0.6139218688180343
$data[33,9] = 0.3244650893419825
$data[33,10] = 0.3144596800720633
$data[34,0] = "2022-04-01 00:00:00_diff"
$data[34,1] = -0.6814298979208291
$data[34,2] = -0.3328666912667296
$data[34,3] = 1.020216564564436
$data[34,4] = 1.335622890042328
$data[34,5] = -0.2643000581021264
$data[34,6] = 0.705045326079295
$data[34,7] = 0.6213407849888899
$data[34,8] = 0.3301586542973803
$data[34,9] = 0.3221046561371852
$data[34,10] = -0.697211019090521
$data[35,0] = "2022-07-01 00:00:00_diff"
$data[35,1] = -0.5026722255049707
$data[35,2] = 0.9306463614384912
$data[35,3] = 1.366285470889476
$data[35,4] = -0.2835604934770959
$data[35,5] = 0.6804660345951945
$data[35,6] = 0.606832781780235
$data[35,7] = 0.3135118745842211
$data[35,8] = 0.3043967211331541
$data[35,9] = -0.7142363390923157
$data[35,10] = -0.6927154267442934
$data[36,0] = "2022-10-01 00:00:00_diff"
$data[36,1] = 0.695268442486379
$data[36,2] = 1.341678812480128
$data[36,3] = -0.2408120245028667
$data[36,4] = 0.662677011776885
$data[36,5] = 0.5962167627018904
$data[36,6] = 0.3117261734516353
$data[36,7] = 0.2984398913470776
$data[36,8] = -0.7204764473137611
$data[36,9] = -0.6980947118431513
$data[36,10] = 0.8406195607711028
$data[37,0] = "2023-01-01 00:00:00_diff"
$data[37,1] = 1.533847576318271
$data[37,2] = -0.1936259769839437
$data[37,3] = 0.583108652925594
$data[37,4] = 0.5852151936452853
$data[37,5] = 0.3035951384809393
$data[37,6] = 0.2760442649778436
$data[37,7] = -0.7388925974841081
$data[37,8] = -0.7150001535555365
$data[37,9] = 0.8224589278409099
$data[37,10] = -0.1294633805524922
$data[38,0] = "2023-04-01 00:00:00_diff"
$data[38,1] = 0.2988169550551557
$data[38,2] = 0.6695032207238313
$data[38,3] = 0.413169930687634
$data[38,4] = 0.2858903600785752
$data[38,5] = 0.2603689453763344
$data[38,6] = -0.7842151956877904
$data[38,7] = -0.7517550993415597
$data[38,8] = 0.7886969707147926
$data[38,9] = -0.165741845517162
$data[38,10] = ""
$data[39,0] = "2023-07-01 00:00:00_diff"
$data[39,1] = 0.7325543163948336
$data[39,2] = 0.4229506677669676
$data[39,3] = 0.2454962336276151
$data[39,4] = 0.2441103205764524
$data[39,5] = -0.7983350006221057
$data[39,6] = -0.7708676270699714
$data[39,7] = 0.7706321548531918
$data[39,8] = -0.1832014377874837
$data[39,9] = ""
$data[39,10] = ""
$data[40,0] = "2023-10-01 00:00:00_diff"
$data[40,1] = 0.6205859453382192
$data[40,2] = 0.2998947699497946
$data[40,3] = 0.1714200530674933
$data[40,4] = -0.8043240939096088
$data[40,5] = -0.7743864887232023
$data[40,6] = 0.7539745674340489
$data[40,7] = -0.1961654531415352
$data[40,8] = ""
$data[40,9] = ""
$data[40,10] = ""
$data[41,0] = "2024-01-01 00:00:00_diff"
$data[41,1] = 0.4819063909728227
$data[41,2] = 0.2168522670283788
$data[41,3] = -0.8642139633931005
$data[41,4] = -0.7761405111915634
$data[41,5] = 0.753463083855529
$data[41,6] = -0.2077890329460543
$data[41,7] = ""
$data[41,8] = ""
$data[41,9] = ""
$data[41,10] = ""
$data[42,0] = "2024-04-01 00:00:00_diff"
$data[42,1] = 0.4586307429802475
$data[42,2] = -0.854661727990704
$data[42,3] = -0.8262811245050392
$data[42,4] = 0.7632630423627871
$data[42,5] = -0.2054539185520808
$data[42,6] = ""
$data[42,7] = ""
$data[42,8] = ""
$data[42,9] = ""
$data[42,10] = ""
$data[43,0] = "2024-07-01 00:00:00_diff"
$data[43,1] = -0.7081246783130091
$data[43,2] = -0.7766873914013612
$data[43,3] = 0.7351734699582756
$data[43,4] = -0.1922981331890715
$data[43,5] = ""
$data[43,6] = ""
$data[43,7] = ""
$data[43,8] = ""
$data[43,9] = ""
$data[43,10] = ""
$data[44,0] = "2024-10-01 00:00:00_diff"
$data[44,1] = -0.9804124188973848
$data[44,2] = 0.7061650899858334
$data[44,3] = -0.1048233766940621
$data[44,4] = ""
$data[44,5] = ""
$data[44,6] = ""
$data[44,7] = ""
$data[44,8] = ""
$data[44,9] = ""
$data[44,10] = ""
$data[45,0] = "2025-01-01 00:00:00_diff"
$data[45,1] = 0.3851484876896414
$data[45,2] = -0.1636275381640872
$data[45,3] = ""
$data[45,4] = ""
$data[45,5] = ""
$data[45,6] = ""
$data[45,7] = ""
$data[45,8] = ""
$data[45,9] = ""
$data[45,10] = ""
$data[46,0] = "2025-04-01 00:00:00_diff"
$data[46,1] = -0.01814927404267447
$data[46,2] = ""
$data[46,3] = ""
$data[46,4] = ""
$data[46,5] = ""
$data[46,6] = ""
$data[46,7] = ""
$data[46,8] = ""
$data[46,9] = ""
$data[46,10] = ""
$data[47,0] = "2025-07-01 00:00:00_diff"
$data[47,1] = ""
$data[47,2] = ""
$data[47,3] = ""
$data[47,4] = ""
$data[47,5] = ""
$data[47,6] = ""
$data[47,7] = ""
$data[47,8] = ""
$data[47,9] = ""
$data[47,10] = ""

$ws.Range("A1:K48").Value = $data
$ws.Range("A49:K53").Clear()
